$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts the existing rows 4..76 down to 5..77
# and keeps formatting consistent with Excel's native "insert row" behaviour.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new market record.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44756
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100114001
$ws.Cells.Item(4, 7).Value = "Papa"
$ws.Cells.Item(4, 8).Value = "Asterix"
$ws.Cells.Item(4, 9).Value = "1a (guarda)"
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 11).Value = 11000
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11500
$ws.Cells.Item(4, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(4, 16).Value = 460
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
